$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 168, shifting existing rows 168:221 down to 169:222
$ws.Rows.Item(168).Insert()

# Populate the newly inserted row 168 with the new weekly record
$ws.Range("A168").Value = 3
$ws.Range("B168").Value = "Femacal de La Calera"
$ws.Range("C168").Value = "Coquimbo"
$ws.Range("D168").Value = 44468
$ws.Range("E168").Value = 5
$ws.Range("F168").Value = 100112003
$ws.Range("G168").Value = "Ajo"
$ws.Range("H168").Value = "Chino"
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 53
$ws.Range("K168").Value = 16500
$ws.Range("L168").Value = 17000
$ws.Range("M168").Value = 16764
$ws.Range("N168").Value = "`$/caja 10 kilos"
$ws.Range("O168").Value = "China"
$ws.Range("P168").Value = 1676
$ws.Range("Q168").Value = 10
$ws.Range("R168").Value = "Hortaliza"
